# Installation.pptx edit:
#  - Slide 1 title: merge the four runs "SAGA " / "Components" / ":" /
#    " Installation and Deployment" into a single run of text.
#  - Slide 2 title: "Agenda" -> "Outline".
#  - Slide 3 title: merge the three runs "Mephisto" / ":" /
#    " The Easy Way Out" into a single run of text.

$p = $ppt.ActivePresentation

# --- Slide 1: title "SAGA Components: Installation and Deployment" ---
$slide1 = $p.Slides.Item(1)
$title1 = $slide1.Shapes.Item(2)
# The visible text is already correct, but it is split across several
# runs. Re-assigning through a distinct intermediate value forces the
# text-range setter to rebuild the paragraph as a single run.
$title1.TextFrame.TextRange.Text = "TEMP"
$title1.TextFrame.TextRange.Text = "SAGA Components: Installation and Deployment"

# --- Slide 2: title "Agenda" -> "Outline" ---
$slide2 = $p.Slides.Item(2)
$title2 = $slide2.Shapes.Item(1)
$title2.TextFrame.TextRange.Text = "Outline"

# --- Slide 3: title "Mephisto: The Easy Way Out" ---
$slide3 = $p.Slides.Item(3)
$title3 = $slide3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "TEMP"
$title3.TextFrame.TextRange.Text = "Mephisto: The Easy Way Out"
